$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Insert two new header rows at the top, shifting all existing data down ---
$ws.Range("A1:A2").EntireRow.Insert()

# --- 2) Populate the two new header rows ---
$ws.Cells.Item(1, 1).Value = "Unnamed: 0"
$ws.Cells.Item(1, 2).Value = "Unnamed: 1"
$ws.Cells.Item(2, 1).Value = "cidades"
$ws.Cells.Item(2, 2).Value = "Casos confirmados"

# --- 3) Style the very first row (A1:B1): bold font, thin box border, centered/top aligned ---
$hdr = $ws.Range("A1:B1")
$hdr.Font.Bold = $true
$hdr.HorizontalAlignment = -4108
$hdr.VerticalAlignment = -4160
$hdr.Borders.LineStyle = 1

# --- 4) Append the new "obitos" (deaths) block below the existing table ---
$tail = @(
    ,@(98, 'outros estados', 34)
    ,@(99, 'outros paises', 35)
    ,@(100, 'cidade', 'Óbtos')
    ,@(101, 'sao paulo', 212)
    ,@(102, 'guarulhos', 5)
    ,@(103, 'campinas', 4)
    ,@(104, 'sao bernardo do campo', 4)
    ,@(105, 'santo andre', 3)
    ,@(106, 'cotia', 2)
    ,@(107, 'osasco', 2)
    ,@(108, 'santos', 2)
    ,@(109, 'sorocaba', 2)
    ,@(110, 'taboao da serra', 2)
    ,@(111, 'americana', 1)
    ,@(112, 'aruja', 1)
    ,@(113, 'barueri', 1)
    ,@(114, 'caieiras', 1)
    ,@(115, 'carapicuiba', 1)
    ,@(116, 'cravinhos', 1)
    ,@(117, 'diadema', 1)
    ,@(118, 'dracena', 1)
    ,@(119, 'embu das artes', 1)
    ,@(120, 'francisco morato', 1)
    ,@(121, 'franco da rocha', 1)
    ,@(122, 'itapecerica da serra', 1)
    ,@(123, 'itapevi', 1)
    ,@(124, 'jaboticabal', 1)
    ,@(125, 'mairipora', 1)
    ,@(126, 'mogi das cruzes', 1)
    ,@(127, 'nova odessa', 1)
    ,@(128, 'penapolis', 1)
    ,@(129, 'ribeirao preto', 1)
    ,@(130, 'sao caetano do sul', 1)
    ,@(131, 'sao sebastiao', 1)
    ,@(132, 'vargem grande paulista', 1)
)

foreach ($row in $tail) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
}

Write-Host "done"
